$wb = $excel.ActiveWorkbook

# The "Status" value moves from "Ready for handoff" to "In Translation"
# for the one tracked file, on every sheet that surfaces it:
#   - Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
#   - zh-cn!C2 (Status column)
#   - de-de!C2 (Status column)
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus

# Shrinking the status text narrows the report's auto-fitted Status
# columns: Overview's per-language columns (E:F) and the Status column
# (C) on each per-language sheet.
$newColumnWidth = 12.5

$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
